# Auto-update draw results: append the 2025-10-03 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 17

# Columns A (date) and C (phase code) look numeric/date-like; force them to
# Text format before assignment so Excel stores the literal strings instead
# of auto-converting to a date serial / number (matches the rest of the
# table, which stores every column as text).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-10-03"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "251003"
$ws.Cells.Item($newRow, 4).Value = "1-9-6"
$ws.Cells.Item($newRow, 5).Value = "2025-10-03T21:36:35.093+04:00"
